$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 (matching style/format of existing header row)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-18
$data = @(
    @(2, 4, 6),
    @(3, 6, 8),
    @(4, 6, 7),
    @(5, 7, 8),
    @(6, 4, 5),
    @(7, 7, 8),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 7, 7),
    @(12, 3, 3),
    @(13, 4, 4),
    @(14, 4, 4),
    @(15, 7, 7),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
